$d = $word.ActiveDocument

# The second paragraph currently holds the text "S" (followed by the
# _GoBack bookmark). We need to:
#   1. Insert a new empty paragraph right before it.
#   2. Replace the "S" text with "Adding rubbish and commiting ."

$target = $d.Paragraphs.Item(2)

# Insert a new paragraph break before the existing text, which pushes
# the current paragraph's content down into a new paragraph and leaves
# an empty paragraph with the same paragraph formatting in its place.
$insertRange = $target.Range
$insertRange.Collapse(1)  # wdCollapseStart
$insertRange.Text = "`r"

# Now replace the "S" run text with the new sentence.
$d.Content.Find.Execute("S", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Adding rubbish and commiting .", 2)
